$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 45.36402809467413
$ws.Range("C2").Value = 21.25658086716114
$ws.Range("D2").Value = 12.02763682269174
$ws.Range("E2").Value = 10.44940240420458
$ws.Range("G2").Value = 3.928581206466307
$ws.Range("J2").Value = 6.791925973068803
$ws.Range("L2").Value = 16.71661688560557
$ws.Range("N2").Value = 28.3175446150848
$ws.Range("B3").Value = 45.11839541230295
$ws.Range("C3").Value = 20.9730771423608
$ws.Range("D3").Value = 12.05119707768959
$ws.Range("E3").Value = 10.46959479310382
$ws.Range("G3").Value = 3.935609978724304
$ws.Range("J3").Value = 6.776634613213098
$ws.Range("L3").Value = 16.73395220772517
$ws.Range("N3").Value = 28.24225038089121
$ws.Range("B4").Value = 44.98030645520984
$ws.Range("C4").Value = 20.80463372797779
$ws.Range("D4").Value = 12.06775729287448
$ws.Range("E4").Value = 10.48269380494959
$ws.Range("G4").Value = 3.940136018303979
$ws.Range("J4").Value = 6.767069876757239
$ws.Range("L4").Value = 16.74801561131698
$ws.Range("N4").Value = 28.19704269424473
$ws.Range("B5").Value = 44.92726309355003
$ws.Range("C5").Value = 20.73748066378642
$ws.Range("D5").Value = 12.07503116076184
$ws.Range("E5").Value = 10.48820865939218
$ws.Range("G5").Value = 3.942033597875967
$ws.Range("J5").Value = 6.76312756286462
$ws.Range("L5").Value = 16.75460363003413
$ws.Range("N5").Value = 28.17887834457014
$ws.Range("B6").Value = 44.91865101897908
$ws.Range("C6").Value = 20.72642195553177
$ws.Range("D6").Value = 12.07627068266781
$ws.Range("E6").Value = 10.48913510105391
$ws.Range("G6").Value = 3.942351910001771
$ws.Range("J6").Value = 6.762470238696587
$ws.Range("L6").Value = 16.75574924507839
$ws.Range("N6").Value = 28.17587778001117
$ws.Range("B7").Value = 44.97957798675992
$ws.Range("C7").Value = 20.80372195494717
$ws.Range("D7").Value = 12.06785326480293
$ws.Range("E7").Value = 10.48276746303454
$ws.Range("G7").Value = 3.940161394011405
$ws.Range("J7").Value = 6.767016890479888
$ws.Range("L7").Value = 16.74810099320417
$ws.Range("N7").Value = 28.19679667728429
$ws.Range("B8").Value = 45.27670996494214
$ws.Range("C8").Value = 21.15770567893496
$ws.Range("D8").Value = 12.03532508679108
$ws.Range("E8").Value = 10.45621969676025
$ws.Range("G8").Value = 3.930961244653715
$ws.Range("J8").Value = 6.786689748535248
$ws.Range("L8").Value = 16.72188278202132
$ws.Range("N8").Value = 28.29137132871922
$ws.Range("B9").Value = 45.95894158207147
$ws.Range("C9").Value = 21.8931416940407
$ws.Range("D9").Value = 11.98820419611294
$ws.Range("E9").Value = 10.40968944344863
$ws.Range("G9").Value = 3.914575554973974
$ws.Range("J9").Value = 6.82389288587176
$ws.Range("L9").Value = 16.69771892908696
$ws.Range("N9").Value = 28.48491372854141
$ws.Range("B10").Value = 46.51860825849026
$ws.Range("C10").Value = 22.45393319780663
$ws.Range("D10").Value = 11.96381451208149
$ws.Range("E10").Value = 10.37883208530067
$ws.Range("G10").Value = 3.903527813202675
$ws.Range("J10").Value = 6.85041633300516
$ws.Range("L10").Value = 16.69673092016818
$ws.Range("N10").Value = 28.63198597924905
$ws.Range("B11").Value = 46.78532944910198
$ws.Range("C11").Value = 22.7124742518848
$ws.Range("D11").Value = 11.95495467413972
$ws.Range("E11").Value = 10.36550788319948
$ws.Range("G11").Value = 3.898713030393746
$ws.Range("J11").Value = 6.862313027791774
$ws.Range("L11").Value = 16.69995028527619
$ws.Range("N11").Value = 28.69993772446768
$ws.Range("B12").Value = 46.88801899593118
$ws.Range("C12").Value = 22.810785299766
$ws.Range("D12").Value = 11.95192218116818
$ws.Range("E12").Value = 10.36056416750896
$ws.Range("G12").Value = 3.896919804492406
$ws.Range("J12").Value = 6.866794172956797
$ws.Range("L12").Value = 16.70169881150214
$ws.Range("N12").Value = 28.72581813323545
$ws.Range("B13").Value = 46.86582881163408
$ws.Range("C13").Value = 22.78959551615713
$ws.Range("D13").Value = 11.9525609204925
$ws.Range("E13").Value = 10.36162436480361
$ws.Range("G13").Value = 3.897304677042175
$ws.Range("J13").Value = 6.865830133271264
$ws.Range("L13").Value = 16.70129866194095
$ws.Range("N13").Value = 28.72023775261239
$ws.Range("B14").Value = 46.7937442804199
$ws.Range("C14").Value = 22.72055469335303
$ws.Range("D14").Value = 11.95469871873977
$ws.Range("E14").Value = 10.36509912283422
$ws.Range("G14").Value = 3.898564900531645
$ws.Range("J14").Value = 6.862682170363188
$ws.Range("L14").Value = 16.70008351741256
$ws.Range("N14").Value = 28.70206395677621
$ws.Range("B15").Value = 46.74980854326887
$ws.Range("C15").Value = 22.67831566634847
$ws.Range("D15").Value = 11.95605021815957
$ws.Range("E15").Value = 10.36724075964476
$ws.Range("G15").Value = 3.899340725437728
$ws.Range("J15").Value = 6.860750855531054
$ws.Range("L15").Value = 16.6994082036051
$ws.Range("N15").Value = 28.69095126252996
$ws.Range("B16").Value = 46.50141672634579
$ws.Range("C16").Value = 22.43709860563342
$ws.Range("D16").Value = 11.96443859724979
$ws.Range("E16").Value = 10.37971714406569
$ws.Range("G16").Value = 3.903846688817989
$ws.Range("J16").Value = 6.849635492063136
$ws.Range("L16").Value = 16.69659451674038
$ws.Range("N16").Value = 28.62756608104915
$ws.Range("B17").Value = 46.35210364533675
$ws.Range("C17").Value = 22.28993738292712
$ws.Range("D17").Value = 11.97015788425233
$ws.Range("E17").Value = 10.38755314410367
$ws.Range("G17").Value = 3.906664754181182
$ws.Range("J17").Value = 6.842773601038562
$ws.Range("L17").Value = 16.69580947552983
$ws.Range("N17").Value = 28.58894825247245
$ws.Range("B18").Value = 46.26736763580112
$ws.Range("C18").Value = 22.20562253212654
$ws.Range("D18").Value = 11.97365778540208
$ws.Range("E18").Value = 10.39212735066336
$ws.Range("G18").Value = 3.908305498790194
$ws.Range("J18").Value = 6.838810894612175
$ws.Range("L18").Value = 16.69570317138247
$ws.Range("N18").Value = 28.56683506574416
$ws.Range("B19").Value = 46.23887576046819
$ws.Range("C19").Value = 22.17713403008882
$ws.Range("D19").Value = 11.97487887929227
$ws.Range("E19").Value = 10.39368765081332
$ws.Range("G19").Value = 3.908864448097009
$ws.Range("J19").Value = 6.837466441899211
$ws.Range("L19").Value = 16.69572641198101
$ws.Range("N19").Value = 28.55936494541422
$ws.Range("B20").Value = 46.36788017479
$ws.Range("C20").Value = 22.3055695899588
$ws.Range("D20").Value = 11.96952728095636
$ws.Range("E20").Value = 10.38671204384207
$ws.Range("G20").Value = 3.906362712207068
$ws.Range("J20").Value = 6.843505706528474
$ws.Range("L20").Value = 16.69585730027201
$ws.Range("N20").Value = 28.59304897904398
$ws.Range("B21").Value = 46.81487190263252
$ws.Range("C21").Value = 22.7408232916382
$ws.Range("D21").Value = 11.95406203330783
$ws.Range("E21").Value = 10.36407574263278
$ws.Range("G21").Value = 3.898193929781805
$ws.Range("J21").Value = 6.863607447961881
$ws.Range("L21").Value = 16.70042605286939
$ws.Range("N21").Value = 28.70739802712765
$ws.Range("B22").Value = 47.1168178046948
$ws.Range("C22").Value = 23.02762071254741
$ws.Range("D22").Value = 11.94583507515951
$ws.Range("E22").Value = 10.34987508820818
$ws.Range("G22").Value = 3.893030043645768
$ws.Range("J22").Value = 6.876606135426762
$ws.Range("L22").Value = 16.70649856646093
$ws.Range("N22").Value = 28.78299653715303
$ws.Range("B23").Value = 46.95478513346802
$ws.Range("C23").Value = 22.87436629161191
$ws.Range("D23").Value = 11.9500535333905
$ws.Range("E23").Value = 10.35740015694204
$ws.Range("G23").Value = 3.895770203856253
$ws.Range("J23").Value = 6.86968106611702
$ws.Range("L23").Value = 16.70297458042356
$ws.Range("N23").Value = 28.74256980774595
$ws.Range("B24").Value = 46.3607441564813
$ws.Range("C24").Value = 22.29850135898714
$ws.Range("D24").Value = 11.96981171706112
$ws.Range("E24").Value = 10.38709208972247
$ws.Range("G24").Value = 3.906499201192775
$ws.Range("J24").Value = 6.84317477684275
$ws.Range("L24").Value = 16.69583460416841
$ws.Range("N24").Value = 28.59119476406769
$ws.Range("B25").Value = 45.76393685467961
$ws.Range("C25").Value = 21.69024670564546
$ws.Range("D25").Value = 11.99916002331836
$ws.Range("E25").Value = 10.42168952531537
$ws.Range("G25").Value = 3.918832970437621
$ws.Range("J25").Value = 6.813972725128665
$ws.Range("L25").Value = 16.70132192599938
$ws.Range("N25").Value = 28.43169451730951
